$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O), Precio promedio ponderado (P), Precio $/Kg (S)
$updates = @(
    @{Row=2;  D=44320; M=20;  N=12000; O=12000; P=12000; S=1714},
    @{Row=3;  D=44320; M=30;  N=8000;  O=8000;  P=8000;  S=1143},
    @{Row=4;  D=44322; M=45;  N=12000; O=12000; P=12000; S=1714},
    @{Row=5;  D=44322; M=80;  N=8000;  O=8000;  P=8000;  S=1143},
    @{Row=11; D=44301; M=100; N=14000; O=14000; P=14000; S=2000},
    @{Row=12; D=44301; M=80;  N=12000; O=12000; P=12000; S=1714},
    @{Row=13; D=44302; M=50;  N=15000; O=15000; P=15000; S=2143},
    @{Row=14; D=44302; M=30;  N=12000; O=12000; P=12000; S=1714},
    @{Row=15; D=44292; M=25;  N=16000; O=16000; P=16000; S=2286},
    @{Row=16; D=44292; M=30;  N=15000; O=15000; P=15000; S=2143},
    @{Row=17; D=44300; M=100; N=15000; O=15000; P=15000; S=2143},
    @{Row=18; D=44300; M=80;  N=12000; O=12000; P=12000; S=1714}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value = $u.D   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $u.M  # M - Volumen
    $ws.Cells.Item($r, 14).Value = $u.N  # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $u.O  # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $u.P  # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $u.S  # S - Precio $/Kg
}
